# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" on the Overview sheet and the
# "Latest Handoff Datetime" on the per-locale sheets for the file that was
# just handed off (4335d125-1ecd-4dc3-bb84-428b79f8d32a.md -> row 7).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-18 02:36:18"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-18 02:36:13"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-18 02:36:18"
